$d = $word.ActiveDocument

$emdash = [char]0x2014
$rsquote = [char]0x2019

# 1. "shadow side — a rise" -> "shadow side; a rise"
$search1 = "shadow side " + $emdash + " a rise"
$replace1 = "shadow side; a rise"
$d.Content.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)

# 2. "already done — accounts" -> "already done; accounts"
$search2 = "already done " + $emdash + " accounts"
$replace2 = "already done; accounts"
$d.Content.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)

# 3. "not static—they" -> "not static, they"
$search3 = "not static" + $emdash + "they"
$replace3 = "not static, they"
$d.Content.Find.Execute($search3, $true, $false, $false, $false, $false, $true, 1, $false, $replace3, 2)

# 4. "reported and analyzed — by " -> "reported and analyzed  by " (em dash removed, double space kept)
$search4 = "reported and analyzed " + $emdash + " by "
$replace4 = "reported and analyzed  by "
$d.Content.Find.Execute($search4, $true, $false, $false, $false, $false, $true, 1, $false, $replace4, 2)

# 5. "analyze SSL" -> "analyse SSL"
$search5 = "inspect a URL, analyze SSL certificates"
$replace5 = "inspect a URL, analyse SSL certificates"
$d.Content.Find.Execute($search5, $true, $false, $false, $false, $false, $true, 1, $false, $replace5, 2)

# 6. "flags phishing attempts — while ... freely — could make" -> em dashes removed
$search6 = "flags phishing attempts " + $emdash + " while allowing the user to browse freely " + $emdash + " could make"
$replace6 = "flags phishing attempts while allowing the user to browse freely could make"
$d.Content.Find.Execute($search6, $true, $false, $false, $false, $false, $true, 1, $false, $replace6, 2)

# 7. "By analyzing the current webpage's" -> "By analysing the current webpage's"
$search7 = "detector. By analyzing the current webpage" + $rsquote + "s structure"
$replace7 = "detector. By analysing the current webpage" + $rsquote + "s structure"
$d.Content.Find.Execute($search7, $true, $false, $false, $false, $false, $true, 1, $false, $replace7, 2)

# 8. "high-end" -> "high end"
$search8 = "access to high-end security tools is limited"
$replace8 = "access to high end security tools is limited"
$d.Content.Find.Execute($search8, $true, $false, $false, $false, $false, $true, 1, $false, $replace8, 2)

# 9. Move the "_GoBack" bookmark from the References section to the end of the
#    "As a result..." paragraph (1.3 Project Motivation), matching Word's
#    behaviour of tracking the last edit location with this hidden bookmark.
$found = $d.Content
$found.Find.Execute("Webpage Phishing Detector System.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($found.End, $found.End)
$d.Bookmarks.Add("_GoBack", $target)
